# "auto column width, extract live value"
#
# The META sheet holds a scraped stock quote; re-running the scrape pulled
# fresher Bid/Ask quotes, so push the newly-extracted live values into the
# sheet, then auto-size the columns to fit the (possibly wider/narrower)
# text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("META")
$ws.Activate()

# Bid / Ask rows (A3:B3 = "Bid", A4:B4 = "Ask") get the freshly extracted
# live values.
$ws.Range("B3").Value = "120.15 x 800"
$ws.Range("B4").Value = "120.42 x 1000"

# Auto-fit the columns now that the values changed.
$ws.Columns.AutoFit()
